$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in task rows 10-14 (tasks #7-#11) with Desc (B), Owner (C), DueDate (D)
$dueDate = Get-Date -Year 2018 -Month 8 -Day 14 -Hour 0 -Minute 0 -Second 0

$ws.Range("B10").Value = "Analysis of tools available (Cost/Adv/Disadv/Reliability etc)"
$ws.Range("C10").Value = "Mihir + Nikhil "
$ws.Range("D10").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D10").Value = $dueDate

$ws.Range("B11").Value = "Compilation of word files for Exhibit F "
$ws.Range("C11").Value = "Mihir"
$ws.Range("D11").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D11").Value = $dueDate

$ws.Range("B12").Value = "Simple Twitter App "
$ws.Range("C12").Value = "Sangeeta"
$ws.Range("D12").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D12").Value = $dueDate

$ws.Range("B13").Value = "Look into Facebook Graph API Access "
$ws.Range("C13").Value = "Mihir "
$ws.Range("D13").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D13").Value = $dueDate

$ws.Range("B14").Value = "GitHub - organize "
$ws.Range("C14").Value = "All "
$ws.Range("D14").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D14").Value = $dueDate

# Update selected cell to match final cursor position
$ws.Range("B15").Select()
